# ASC-232: make the internship duration read "{{ duration }}months", and
# flip the signing order in the signature block so the employer
# representative's name precedes the intern's name.

$d = $word.ActiveDocument

# --- 1) "{{ duration }}" -> "{{ duration }}months" -------------------------
# Insert the new "months" text right after the existing placeholder instead
# of rewriting the whole phrase, so the pre-existing runs stay untouched.
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -like "*Duration:*" -and $txt -like "*{{ duration }}*") {
        $pText = $p.Range.Text
        $token = "{{ duration }}"
        $offset = $pText.IndexOf($token)
        $tokenRange = $d.Range($p.Range.Start + $offset, $p.Range.Start + $offset + $token.Length)
        $tokenRange.InsertAfter("months")
        break
    }
}

# --- 2) Swap the "{{ intern_name }}" / "{{ employer_representative_name }}"
#        placeholders in the signature-line paragraph (the one holding both
#        tokens), leaving the tabs/spacing that separates them untouched.
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -like "*intern_name*" -and $txt -like "*employer_representative_name*") {
        $r1 = $p.Range
        $null = $r1.Find.Execute(
            "{{ intern_name }}", $true, $false, $false, $false, $false,
            $true, 1, $false, "%%SWAP_PLACEHOLDER%%", 2)

        $r2 = $p.Range
        $null = $r2.Find.Execute(
            "{{ employer_representative_name }}", $true, $false, $false, $false, $false,
            $true, 1, $false, "{{ intern_name }}", 2)

        $r3 = $p.Range
        $null = $r3.Find.Execute(
            "%%SWAP_PLACEHOLDER%%", $true, $false, $false, $false, $false,
            $true, 1, $false, "{{ employer_representative_name }}", 2)

        break
    }
}
